$wb = $excel.ActiveWorkbook

$oldGuid = "8ef74566-ae50-4db2-98df-a520d80ebfb5"
$newGuid = "4db14c90-50ac-469d-8c12-c56368b6f730"
$oldHash = "699b65ceeb6b22927dfefb3423fd82f0022cc1ba"
$newHash = "63b8aa167639b1a62a163163dc197b223bc41d51"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9228812ae35f564b7f6aed3ce65221a7a220ca1a/e2e/$oldGuid.md"

# --- Overview sheet ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Cells.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md")
$wsOverview.Range("G2").Value = "2016-09-07 03:13:53"

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Cells.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, "$newGuid.md")
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-07 03:13:49"

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Cells.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, "$newGuid.md")
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-07 03:13:53"
